# Fruta / hortaliza, semanal
# Insert two new daily price records (rows 240 and 241) for "Ají" in the
# "Macroferia Regional de Talca" sheet. All existing data from row 240
# downward shifts down by two rows (to 242..328); the two freshly
# inserted rows are then populated with the new observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 240 and below down by two rows, preserving formatting
# (e.g. the date number format carried on column D).
$ws.Rows.Item(240).Resize(2, 1).Insert()

# --- New row 240: Ají, "Cacho cabra verde" -------------------------------
$ws.Cells.Item(240, 1).Value  = 5
$ws.Cells.Item(240, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(240, 3).Value  = "Maule"
$ws.Cells.Item(240, 4).Value  = 44988
$ws.Cells.Item(240, 5).Value  = 7
$ws.Cells.Item(240, 6).Value  = 100112021
$ws.Cells.Item(240, 7).Value  = "Ají"
$ws.Cells.Item(240, 8).Value  = "Cacho cabra verde"
$ws.Cells.Item(240, 9).Value  = "Primera"
$ws.Cells.Item(240, 10).Value = 150
$ws.Cells.Item(240, 11).Value = 12000
$ws.Cells.Item(240, 12).Value = 12000
$ws.Cells.Item(240, 13).Value = 12000
$ws.Cells.Item(240, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(240, 15).Value = "Región del Maule"
$ws.Cells.Item(240, 16).Value = 480
$ws.Cells.Item(240, 17).Value = 25
$ws.Cells.Item(240, 18).Value = "Hortaliza"

# --- New row 241: Ají, "Cristal" -----------------------------------------
$ws.Cells.Item(241, 1).Value  = 5
$ws.Cells.Item(241, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(241, 3).Value  = "Maule"
$ws.Cells.Item(241, 4).Value  = 44988
$ws.Cells.Item(241, 5).Value  = 7
$ws.Cells.Item(241, 6).Value  = 100112021
$ws.Cells.Item(241, 7).Value  = "Ají"
$ws.Cells.Item(241, 8).Value  = "Cristal"
$ws.Cells.Item(241, 9).Value  = "Primera"
$ws.Cells.Item(241, 10).Value = 150
$ws.Cells.Item(241, 11).Value = 12000
$ws.Cells.Item(241, 12).Value = 12000
$ws.Cells.Item(241, 13).Value = 12000
$ws.Cells.Item(241, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(241, 15).Value = "Región del Maule"
$ws.Cells.Item(241, 16).Value = 480
$ws.Cells.Item(241, 17).Value = 25
$ws.Cells.Item(241, 18).Value = "Hortaliza"
